$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string must be forced to
# Text format first, matching the source data which stores every Price/Volume
# cell as text (inlineStr) even when it looks numeric (e.g. "1.00", "683.83").
function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "69.327.65"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").Value = "3.684.58"
$ws.Range("E3").Value = "  -3.34%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue "D5" "683.83"
$ws.Range("E5").Value = "  -3.17%  "
Set-TextValue "D6" "162.24"
$ws.Range("E6").Value = "  -5.35%  "
$ws.Range("D7").Value = "3.683.48"
$ws.Range("E7").Value = "  -3.33%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -4.41%  "
$ws.Range("E10").Value = "  -8.01%  "
$ws.Range("E11").Value = "  -5.42%  "
$ws.Range("E12").Value = "  -3.07%  "
$ws.Range("E13").Value = "  -4.81%  "
$ws.Range("E14").Value = "  -6.15%  "
$ws.Range("D15").Value = "4.306.44"
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("D16").Value = "3.682.82"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("D17").Value = "69.430.98"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("E19").Value = "  -6.20%  "
Set-TextValue "D20" "6.61"
$ws.Range("E20").Value = "  -7.43%  "
Set-TextValue "D21" "481.96"
$ws.Range("E21").Value = "  -4.14%  "
Set-TextValue "D22" "9.92"
$ws.Range("E22").Value = "  -7.59%  "
$ws.Range("E23").Value = "  -7.86%  "
$ws.Range("E24").Value = "  -4.73%  "
$ws.Range("D25").Value = "3.830.17"
$ws.Range("E25").Value = "  -3.35%  "
$ws.Range("E26").Value = "  -9.60%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D27" "1.00"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D28" "11.52"
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("E29").Value = "  -8.41%  "
$ws.Range("E30").Value = "  -9.82%  "
$ws.Range("E31").Value = "  -10.78%  "
$ws.Range("E32").Value = "  -7.58%  "
$ws.Range("E33").Value = "  -6.86%  "
$ws.Range("E34").Value = "  -2.58%  "
Set-TextValue "D35" "27.11"
$ws.Range("E35").Value = "  -6.61%  "
Set-TextValue "D36" "0.999"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "3.655.21"
$ws.Range("E37").Value = "  -3.24%  "
Set-TextValue "D38" "8.50"
$ws.Range("E38").Value = "  -7.22%  "
Set-TextValue "D39" "6.35"
$ws.Range("E39").Value = "  +7.09%  "
$ws.Range("E40").Value = "  -1.85%  "
Set-TextValue "D41" "0.0937"
$ws.Range("E41").Value = "  -7.42%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -6.68%  "
Set-TextValue "D45" "161.88"
$ws.Range("E45").Value = "  -3.05%  "
Set-TextValue "D46" "48.32"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D47" "2.83"
$ws.Range("E47").Value = "  -13.34%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D48" "29.99"
$ws.Range("E48").Value = "  +4.96%  "
Set-TextValue "D49" "0.000287"
$ws.Range("E49").Value = "  -8.23%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  -3.07%  "
